$d = $word.ActiveDocument

# 1) "Agile, TDD, MVC, " -> "Agile, TDD, MVC and other patterns, "
$d.Content.Find.Execute(
    "Agile, TDD, MVC, ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Agile, TDD, MVC and other patterns, ", 2
)

# 2) "Eclipse, Visual C/C++, gcc, xc, Ant, make. SubVersion, git, VSS." ->
#    "Eclipse, Visual C/C++, gcc, apache ant, make. SubVersion, git."
$d.Content.Find.Execute(
    "Eclipse, Visual C/C++, gcc, xc, Ant, make. SubVersion, git, VSS.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Eclipse, Visual C/C++, gcc, apache ant, make. SubVersion, git.", 2
)

# 3) " Installsheild, Firebug, cygwin." -> " Installsheild, Firebug, cygwin jQuery."
$d.Content.Find.Execute(
    " Installsheild, Firebug, cygwin.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Installsheild, Firebug, cygwin jQuery.", 2
)

# 4) " Apache, ISS or IBM httpd. " -> " Apache and ISS. "
$d.Content.Find.Execute(
    " Apache, ISS or IBM httpd. ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Apache and ISS. ", 2
)

# 5) "Transactions between web server and " -> "Transaction protocol between web server and "
$d.Content.Find.Execute(
    "Transactions between web server and ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Transaction protocol between web server and ", 2
)

# 6) Remove the "Neil William Hancock" / "Address: ..." / "Mobile: ..." bullet paragraphs
#    (the CV author's name, postal address and mobile number), leaving
#    "Contact Details" and "Email: ..." paragraphs intact.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Neil William Hancock") {
        $startPara = $i
    }
    if ($t -match "Mobile: 07532 242296") {
        $endPara = $i
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeToRemove = $d.Range($d.Paragraphs.Item($startPara).Range.Start, $d.Paragraphs.Item($endPara).Range.End)
    $rangeToRemove.Delete()
}
